$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''42.884.04'
$ws.Range("E2").Value = '''  -1.60%  '
$ws.Range("D3").Value = '''2.339.16'
$ws.Range("E3").Value = '''  +1.02%  '
$ws.Range("E4").Value = '''  -0.07%  '
$ws.Range("D5").Value = '''306.85'
$ws.Range("E5").Value = '''  -1.60%  '
$ws.Range("D6").Value = '''100.45'
$ws.Range("E6").Value = '''  -2.55%  '
$ws.Range("E7").Value = '''  -5.27%  '
$ws.Range("D9").Value = '''0.512'
$ws.Range("E9").Value = '''  -4.25%  '
$ws.Range("D10").Value = '''35.00'
$ws.Range("E10").Value = '''  -2.89%  '
$ws.Range("D11").Value = '''52.14'
$ws.Range("E11").Value = '''  +0.31%  '
$ws.Range("D12").Value = '''0.0799'
$ws.Range("E12").Value = '''  -2.27%  '
$ws.Range("E13").Value = '''  -0.55%  '
$ws.Range("E14").Value = '''  -3.28%  '
$ws.Range("D15").Value = '''15.84'
$ws.Range("E15").Value = '''  +5.31%  '
$ws.Range("D16").Value = '''2.330.70'
$ws.Range("E16").Value = '''  +0.37%  '
$ws.Range("E17").Value = '''  -1.48%  '
$ws.Range("D18").Value = '''42.793.80'
$ws.Range("E18").Value = '''  -1.62%  '
$ws.Range("D19").Value = '''6.24'
$ws.Range("E19").Value = '''  +1.22%  '
$ws.Range("E20").Value = '''  -2.39%  '
$ws.Range("D21").Value = '''11.69'
$ws.Range("E21").Value = '''  -6.63%  '
$ws.Range("D22").Value = '''67.90'
$ws.Range("E22").Value = '''  -0.63%  '
$ws.Range("D23").Value = '''236.86'
$ws.Range("E23").Value = '''  -2.39%  '
$ws.Range("E24").Value = '''  -2.61%  '
$ws.Range("D25").Value = '''2.57'
$ws.Range("E25").Value = '''  -2.35%  '
$ws.Range("E26").Value = '''  -0.11%  '
$ws.Range("D27").Value = '''25.52'
$ws.Range("E27").Value = '''  +2.92%  '
$ws.Range("E28").Value = '''  +6.42%  '
$ws.Range("D29").Value = '''34.94'
$ws.Range("E29").Value = '''  -6.49%  '
$ws.Range("E30").Value = '''  -3.01%  '
$ws.Range("D31").Value = '''159.93'
$ws.Range("E31").Value = '''  -4.76%  '
$ws.Range("E33").Value = '''  -3.60%  '
$ws.Range("D34").Value = '''4.65'
$ws.Range("E34").Value = '''  +6.39%  '
$ws.Range("B35").Value = '''Hedera'
$ws.Range("C35").Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '''0.0728'
$ws.Range("E35").Value = '''  -2.65%  '
$ws.Range("B36").Value = '''Celestia'
$ws.Range("C36").Value = '''https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D36").Value = '''17.36'
$ws.Range("E36").Value = '''  -1.42%  '
$ws.Range("E37").Value = '''  -3.37%  '
$ws.Range("D38").Value = '''2.97'
$ws.Range("E38").Value = '''  -4.99%  '
$ws.Range("E39").Value = '''  -0.56%  '
$ws.Range("E40").Value = '''  -3.31%  '
$ws.Range("E41").Value = '''  -3.08%  '
$ws.Range("D42").Value = '''2.36'
$ws.Range("E42").Value = '''  +1.86%  '
$ws.Range("D43").Value = '''2.028.69'
$ws.Range("E43").Value = '''  +2.80%  '
$ws.Range("E44").Value = '''  -1.88%  '
$ws.Range("D45").Value = '''18.78'
$ws.Range("E45").Value = '''  -3.78%  '
$ws.Range("D46").Value = '''10.30'
$ws.Range("E46").Value = '''  +3.93%  '
$ws.Range("D47").Value = '''2.94'
$ws.Range("E47").Value = '''  -2.30%  '
$ws.Range("D48").Value = '''56.26'
$ws.Range("E48").Value = '''  +0.99%  '
$ws.Range("E49").Value = '''  -1.03%  '
$ws.Range("D50").Value = '''2.564.41'
$ws.Range("E50").Value = '''  +0.77%  '
$ws.Range("D51").Value = '''4.65'
$ws.Range("E51").Value = '''  +1.40%  '
